$d = $word.ActiveDocument

# 1. Expand the "git branch --lists all the branches" bullet text.
#    Locate the existing text, collapse to its end, and insert the
#    additional wording there so the original run (and its rPr) is
#    preserved instead of being replaced wholesale.
$rng = $d.Content
$find = $rng.Find
$found = $find.Execute(
    "git branch  --lists all the branches",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0
)
$rng.Collapse(0)
$rng.InsertAfter(" the * prefixed the branch-name denotes the current branch  which you have checked out.")

# 2. Add a new bullet "git fetch origin" after the "git mergetool" bullet
#    (which is the last paragraph of the document).
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$newLast = $d.Paragraphs.Last
$newLast.Range.Text = "git fetch origin"
